$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Play Station 5"
$ws.Range("C3").Value = "Electrónica"
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = 1500

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Martillo"
$ws.Range("C4").Value = "Ferretería"
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 1500
